$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-44 (row 43 unchanged at 0)
$newK = @{
    2 = 0
    3 = 1
    4 = 0
    5 = 2
    6 = 0
    7 = 1
    8 = 0
    9 = 2
    10 = 1
    11 = 3
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 2
    36 = 2
    37 = 0
    38 = 1
    39 = 0
    40 = 2
    41 = 1
    42 = 1
    44 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
